$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook already has rows 3..16 filled with scheme names (HKL averaging
# schemes) in column B, an index in column A, and 1's across C:P.
# This run added 3 new schemes ("Spiral-...") and relocated "Gaussian-
# Quadrature" so the block that used to be rows 10..16 becomes rows 10..19,
# in this final order:
$labels = @(
    "Gaussian-Quadrature",
    "Spiral-90deg-10rot-5space",
    "Spiral-90deg-15rot-5space",
    "Spiral-90deg-10rot-3space",
    "NoRotation-tilt60deg",
    "Rotation-NoTilt",
    "Rotation-60detTilt",
    "HexGrid-90degTilt5degRes",
    "HexGrid-90degTilt22p5degRes",
    "HexGrid-60degTilt5degRes"
)

$startRow = 10
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $startRow + $i
    $index = $row - 2

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $index
    # Match the formatting already used on A2:A16 (bold, boxed, centered).
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    $ws.Cells.Item($row, 2).Value = $labels[$i]

    for ($col = 3; $col -le 16; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}
